$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Jude Coal Mine, China, M2011, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 7; $row++) {
    $wsData.Cells.Item($row, 19).Value = $newVersion  # column S = build_version
}
